$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -19.24780237706063
$ws.Range("C2").Value = 2.695720114541946
$ws.Range("D2").Value = -19.24780237706063
$ws.Range("E2").Value = -19.24780237706063
$ws.Range("F2").Value = -19.24780237706063
$ws.Range("G2").Value = -19.24780237706063
$ws.Range("H2").Value = -19.24780237706063
$ws.Range("I2").Value = -19.24780237706063
$ws.Range("J2").Value = -19.24780237706063
$ws.Range("K2").Value = -19.24780237706063

$ws.Range("B3").Value = -19.24780237706063
$ws.Range("C3").Value = -19.24780237706063
$ws.Range("D3").Value = -19.24780237706063
$ws.Range("E3").Value = -19.24780237706063
$ws.Range("F3").Value = -19.24780237706063
$ws.Range("G3").Value = -19.24780237706063
$ws.Range("H3").Value = -19.24780237706063
$ws.Range("I3").Value = 2.405664207703066
$ws.Range("J3").Value = -19.24780237706063
$ws.Range("K3").Value = -19.24780237706063

$ws.Range("B4").Value = -19.24780237706063
$ws.Range("C4").Value = 2.271016435170942
$ws.Range("D4").Value = -19.24780237706063
$ws.Range("E4").Value = -19.24780237706063
$ws.Range("F4").Value = 2.650772810363874
$ws.Range("G4").Value = -19.24780237706063
$ws.Range("H4").Value = 1.911581803367664
$ws.Range("I4").Value = -19.24780237706063
$ws.Range("J4").Value = 2.363325560405597
$ws.Range("K4").Value = -19.24780237706063

$ws.Range("B5").Value = -19.24780237706063
$ws.Range("C5").Value = 1.088474254083443
$ws.Range("D5").Value = -19.24780237706063
$ws.Range("E5").Value = -19.24780237706063
$ws.Range("F5").Value = -19.24780237706063
$ws.Range("G5").Value = 2.101319540487576
$ws.Range("H5").Value = -19.24780237706063
$ws.Range("I5").Value = -19.24780237706063
$ws.Range("J5").Value = -19.24780237706063
$ws.Range("K5").Value = -19.24780237706063

$ws.Range("B6").Value = -19.24780237706063
$ws.Range("C6").Value = -19.24780237706063
$ws.Range("D6").Value = -19.24780237706063
$ws.Range("E6").Value = -19.24780237706063
$ws.Range("F6").Value = -19.24780237706063
$ws.Range("G6").Value = -19.24780237706063
$ws.Range("H6").Value = -19.24780237706063
$ws.Range("I6").Value = -19.24780237706063
$ws.Range("J6").Value = -19.24780237706063
$ws.Range("K6").Value = -19.24780237706063

$ws.Range("B7").Value = 2.906997502922216
$ws.Range("C7").Value = -19.24780237706063
$ws.Range("D7").Value = -19.24780237706063
$ws.Range("E7").Value = -19.24780237706063
$ws.Range("F7").Value = -19.24780237706063
$ws.Range("G7").Value = -19.24780237706063
$ws.Range("H7").Value = -19.24780237706063
$ws.Range("I7").Value = -19.24780237706063
$ws.Range("J7").Value = -19.24780237706063
$ws.Range("K7").Value = -19.24780237706063

$ws.Range("B8").Value = -19.24780237706063
$ws.Range("C8").Value = -19.24780237706063
$ws.Range("D8").Value = -19.24780237706063
$ws.Range("E8").Value = 3.028220051102821
$ws.Range("F8").Value = -19.24780237706063
$ws.Range("G8").Value = -19.24780237706063
$ws.Range("H8").Value = -19.24780237706063
$ws.Range("I8").Value = -19.24780237706063
$ws.Range("J8").Value = -19.24780237706063
$ws.Range("K8").Value = -19.24780237706063

$ws.Range("B9").Value = 3.643788704316509
$ws.Range("C9").Value = -19.24780237706063
$ws.Range("D9").Value = -19.24780237706063
$ws.Range("E9").Value = -19.24780237706063
$ws.Range("F9").Value = -19.24780237706063
$ws.Range("G9").Value = -19.24780237706063
$ws.Range("H9").Value = -19.24780237706063
$ws.Range("I9").Value = -19.24780237706063
$ws.Range("J9").Value = -19.24780237706063
$ws.Range("K9").Value = -19.24780237706063

$ws.Range("B10").Value = -19.24780237706063
$ws.Range("C10").Value = -19.24780237706063
$ws.Range("D10").Value = -19.24780237706063
$ws.Range("E10").Value = -19.24780237706063
$ws.Range("F10").Value = -19.24780237706063
$ws.Range("G10").Value = -19.24780237706063
$ws.Range("H10").Value = -19.24780237706063
$ws.Range("I10").Value = 1.563425697474031
$ws.Range("J10").Value = -19.24780237706063
$ws.Range("K10").Value = 2.234834751239915

$ws.Range("B11").Value = -19.24780237706063
$ws.Range("C11").Value = -19.24780237706063
$ws.Range("D11").Value = -19.24780237706063
$ws.Range("E11").Value = 1.945454055645123
$ws.Range("F11").Value = -19.24780237706063
$ws.Range("G11").Value = 2.638431907167342
$ws.Range("H11").Value = -19.24780237706063
$ws.Range("I11").Value = -19.24780237706063
$ws.Range("J11").Value = -19.24780237706063
$ws.Range("K11").Value = 1.437679599931719

$ws.Range("B12").Value = -19.24780237706063
$ws.Range("C12").Value = -19.24780237706063
$ws.Range("D12").Value = -19.24780237706063
$ws.Range("E12").Value = -19.24780237706063
$ws.Range("F12").Value = -19.24780237706063
$ws.Range("G12").Value = -19.24780237706063
$ws.Range("H12").Value = -19.24780237706063
$ws.Range("I12").Value = -19.24780237706063
$ws.Range("J12").Value = -19.24780237706063
$ws.Range("K12").Value = -19.24780237706063

$ws.Range("B13").Value = -19.24780237706063
$ws.Range("C13").Value = -19.24780237706063
$ws.Range("D13").Value = -19.24780237706063
$ws.Range("E13").Value = 1.844480542989419
$ws.Range("F13").Value = -19.24780237706063
$ws.Range("G13").Value = -19.24780237706063
$ws.Range("H13").Value = -19.24780237706063
$ws.Range("I13").Value = -19.24780237706063
$ws.Range("J13").Value = 2.255695846053625
$ws.Range("K13").Value = 1.637739845047639

$ws.Range("B14").Value = -19.24780237706063
$ws.Range("C14").Value = -19.24780237706063
$ws.Range("D14").Value = -19.24780237706063
$ws.Range("E14").Value = -19.24780237706063
$ws.Range("F14").Value = -19.24780237706063
$ws.Range("G14").Value = -19.24780237706063
$ws.Range("H14").Value = -19.24780237706063
$ws.Range("I14").Value = -19.24780237706063
$ws.Range("J14").Value = -19.24780237706063
$ws.Range("K14").Value = 2.12326220749583

$ws.Range("B15").Value = -19.24780237706063
$ws.Range("C15").Value = -19.24780237706063
$ws.Range("D15").Value = -19.24780237706063
$ws.Range("E15").Value = -19.24780237706063
$ws.Range("F15").Value = -19.24780237706063
$ws.Range("G15").Value = -19.24780237706063
$ws.Range("H15").Value = -19.24780237706063
$ws.Range("I15").Value = -19.24780237706063
$ws.Range("J15").Value = -19.24780237706063
$ws.Range("K15").Value = -19.24780237706063

$ws.Range("B16").Value = -19.24780237706063
$ws.Range("C16").Value = -19.24780237706063
$ws.Range("D16").Value = -19.24780237706063
$ws.Range("E16").Value = -19.24780237706063
$ws.Range("F16").Value = -19.24780237706063
$ws.Range("G16").Value = -19.24780237706063
$ws.Range("H16").Value = -19.24780237706063
$ws.Range("I16").Value = -19.24780237706063
$ws.Range("J16").Value = 2.287752565269422
$ws.Range("K16").Value = -19.24780237706063

$ws.Range("B17").Value = -19.24780237706063
$ws.Range("C17").Value = 1.362452013813138
$ws.Range("D17").Value = -19.24780237706063
$ws.Range("E17").Value = -19.24780237706063
$ws.Range("F17").Value = -19.24780237706063
$ws.Range("G17").Value = -19.24780237706063
$ws.Range("H17").Value = 0.6329963062452429
$ws.Range("I17").Value = 1.075688949621087
$ws.Range("J17").Value = 1.301178468595614
$ws.Range("K17").Value = -19.24780237706063

$ws.Range("B18").Value = -19.24780237706063
$ws.Range("C18").Value = -19.24780237706063
$ws.Range("D18").Value = -19.24780237706063
$ws.Range("E18").Value = -19.24780237706063
$ws.Range("F18").Value = -19.24780237706063
$ws.Range("G18").Value = -19.24780237706063
$ws.Range("H18").Value = 0.5020948677125061
$ws.Range("I18").Value = 1.102853058411805
$ws.Range("J18").Value = 1.449709078721585
$ws.Range("K18").Value = -19.24780237706063

$ws.Range("B19").Value = -19.24780237706063
$ws.Range("C19").Value = -19.24780237706063
$ws.Range("D19").Value = -19.24780237706063
$ws.Range("E19").Value = -19.24780237706063
$ws.Range("F19").Value = -19.24780237706063
$ws.Range("G19").Value = -19.24780237706063
$ws.Range("H19").Value = 1.881254823733245
$ws.Range("I19").Value = 2.022322390471205
$ws.Range("J19").Value = -19.24780237706063
$ws.Range("K19").Value = -19.24780237706063

$ws.Range("B20").Value = -19.24780237706063
$ws.Range("C20").Value = 0.8405782762715375
$ws.Range("D20").Value = 4.3219258933168
$ws.Range("E20").Value = -19.24780237706063
$ws.Range("F20").Value = 3.778201759186945
$ws.Range("G20").Value = -19.24780237706063
$ws.Range("H20").Value = 2.145672401332554
$ws.Range("I20").Value = 1.777259437721976
$ws.Range("J20").Value = -19.24780237706063
$ws.Range("K20").Value = 2.354869759725608

$ws.Range("B21").Value = -19.24780237706063
$ws.Range("C21").Value = 1.141498102205892
$ws.Range("D21").Value = -19.24780237706063
$ws.Range("E21").Value = 2.137215816374193
$ws.Range("F21").Value = -19.24780237706063
$ws.Range("G21").Value = 3.245253881962646
$ws.Range("H21").Value = 2.367867218081626
$ws.Range("I21").Value = -19.24780237706063
$ws.Range("J21").Value = -19.24780237706063
$ws.Range("K21").Value = -19.24780237706063
